$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels (remove spaces)
$ws.Range("C1").Value = "JsonSurvey"
$ws.Range("D1").Value = "JsonStream"
$ws.Range("E1").Value = "StartIndex"
$ws.Range("F1").Value = "DoubleBattery"

# Fix data value in row 7 (Start Index for pt3 right hemisphere)
$ws.Range("E7").Value = 7

# Adjust column widths for E and F (values chosen so the engine's
# pixel-snapped stored width lands as close as possible to the target)
$ws.Range("E1").ColumnWidth = 9.8
$ws.Range("F1").ColumnWidth = 14.6

# Update selected cell to match the saved view state
[void]$ws.Range("E8").Select()
